$d = $word.ActiveDocument

# --- Hunk 1: remove the _GoBack bookmark from the place_of_contract paragraph ---
$t1 = $d.Tables.Item(1)
$cell1 = $t1.Cell(1, 1)
$para1 = $cell1.Range.Paragraphs.Item(1)
$rng1 = $para1.Range
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="002D55CA" w:rsidRDefault="002D55CA" w:rsidP="00A11F87"><w:pPr><w:pStyle w:val="ConsPlusNormal"/></w:pPr><w:r><w:t>г. </w:t></w:r><w:r w:rsidR="00014EDF" w:rsidRPr="00573C02"><w:t>$</w:t></w:r><w:r w:rsidR="00014EDF" w:rsidRPr="00573C02"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00A11F87"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>place_of_cont</w:t></w:r><w:r w:rsidR="00C715A1"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>r</w:t></w:r><w:r w:rsidR="00A11F87"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>act</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00014EDF" w:rsidRPr="00573C02"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($xml1)

# --- Hunk 2: add lang rPr to firstside_requisites paragraph, and append two new paragraphs
#     (one empty with lang rPr, one with the relocated _GoBack bookmark) ---
$t2 = $d.Tables.Item(2)
$cell2 = $t2.Cell(2, 1)
$para2 = $cell2.Range.Paragraphs.Item(1)
$rng2 = $para2.Range
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="002D55CA" w:rsidRPr="008546C7" w:rsidRDefault="008546C7" w:rsidP="00984138"><w:pPr><w:pStyle w:val="ConsPlusNormal"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="008546C7"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>firstside_requisites</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ConsPlusNormal"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ConsPlusNormal"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng2.InsertXML($xml2)
